# LMS-2523 Update BaSynthec Validation
# Rename strain codes on the "openbis-data" sheet (prefix them with "JJS-")
# and add two new strain rows ("MS" and "WT 168 trp+") that duplicate the
# existing OD600 measurement series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")

# Rename the existing strain labels in column A.
$ws.Range("A2").Value = "JJS-MGP1"
$ws.Range("A3").Value = "JJS-MGP100"
$ws.Range("A4").Value = "JJS-MGP20"
$ws.Range("A5").Value = "JJS-MGP999"
$ws.Range("A6").Value = "JJS-MGP1"

# Duplicate the last data row (OD600 series) for the two new strains.
$ws.Range("A6:U6").Copy()
$ws.Range("A7:U7").PasteSpecial()
$ws.Range("A6:U6").Copy()
$ws.Range("A8:U8").PasteSpecial()

$ws.Range("A7").Value = "MS"
$ws.Range("A8").Value = "WT 168 trp+"

# Move the selection/active cell as left by the editor.
$ws.Activate()
[void]$ws.Range("A12").Select()
